$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.785.50"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "2.300.33"
$ws.Range("E3").Value = "  -5.02%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'548.77"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'131.24"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Value = "2.297.28"
$ws.Range("E9").Value = "  -5.06%  "
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").Value = "'5.56"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "'0.336"
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").Value = "'23.83"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").Value = "2.703.35"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").Value = "58.727.08"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").Value = "2.298.10"
$ws.Range("E18").Value = "  -5.20%  "
$ws.Range("D19").Value = "'10.67"
$ws.Range("E19").Value = "  -5.29%  "
$ws.Range("E20").Value = "  -4.25%  "
$ws.Range("D21").Value = "'316.13"
$ws.Range("E21").Value = "  -3.63%  "
$ws.Range("E22").Value = "  -4.26%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'63.08"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").Value = "'0.172"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'8.11"
$ws.Range("E27").Value = "  -7.33%  "
$ws.Range("D29").Value = "'1.75"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").Value = "'169.60"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "0.0₃0731"
$ws.Range("E31").Value = "  -5.58%  "
$ws.Range("D32").Value = "'5.80"
$ws.Range("E32").Value = "  -4.95%  "
$ws.Range("D33").Value = "'1.08"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "'0.385"
$ws.Range("E34").Value = "  -4.94%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'17.87"
$ws.Range("E36").Value = "  -3.67%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -5.65%  "
$ws.Range("E39").Value = "  -5.90%  "
$ws.Range("D40").Value = "'37.98"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("D42").Value = "'298.27"
$ws.Range("E42").Value = "  -9.20%  "
$ws.Range("D43").Value = "'140.72"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("E44").Value = "  -5.69%  "
$ws.Range("D45").Value = "'0.0953"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "'0.0502"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").Value = "'18.61"
$ws.Range("E47").Value = "  -7.19%  "
$ws.Range("D48").Value = "'0.556"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").Value = "'0.0216"
$ws.Range("E49").Value = "  -3.37%  "
$ws.Range("D50").Value = "'16.71"
$ws.Range("E50").Value = "  -4.65%  "
$ws.Range("D51").Value = "'11.02"
$ws.Range("E51").Value = "  -0.22%  "
